# Auto-generated edit script applying the diff changes to Malboro_Profits workbook
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H131").Value = 7216.8
$ws.Range("J131").Value = 7996.3335
$ws.Range("L131").Value = 23989.0005
$ws.Range("N131").Value = -34069.00049999999

$ws.Range("H138").Value = 1842.1555
$ws.Range("I138").Value = 1358.6364
$ws.Range("J138").Value = 3171.8333
$ws.Range("K138").Value = 4075.9092
$ws.Range("L138").Value = 9515.499899999999
$ws.Range("M138").Value = 1064.0908
$ws.Range("N138").Value = -19795.4999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 26769
$ws.Range("I61").Value = 23521.8
$ws.Range("K61").Value = 23521.8
$ws.Range("M61").Value = -23309.8

$ws.Range("H74").Value = 21776.143
$ws.Range("I74").Value = 2552.0833
$ws.Range("J74").Value = 47408.223
$ws.Range("K74").Value = 2552.0833
$ws.Range("L74").Value = 47408.223
$ws.Range("M74").Value = -1678.0833
$ws.Range("N74").Value = -49156.223

$ws.Range("H77").Value = 21776.143
$ws.Range("I77").Value = 2552.0833
$ws.Range("J77").Value = 47408.223
$ws.Range("K77").Value = 12760.4165
$ws.Range("L77").Value = 237041.115
$ws.Range("M77").Value = -8392.416499999999
$ws.Range("N77").Value = -245777.115

$ws.Range("H88").Value = 1839.5555
$ws.Range("I88").Value = 1179.5714
$ws.Range("J88").Value = 2259.5454
$ws.Range("K88").Value = 1179.5714
$ws.Range("L88").Value = 2259.5454
$ws.Range("M88").Value = -773.5714
$ws.Range("N88").Value = -3071.5454

$ws.Range("H91").Value = 1839.5555
$ws.Range("I91").Value = 1179.5714
$ws.Range("J91").Value = 2259.5454
$ws.Range("K91").Value = 1179.5714
$ws.Range("L91").Value = 2259.5454
$ws.Range("M91").Value = 224.4286
$ws.Range("N91").Value = -5067.5454

$ws.Range("H102").Value = 10326.523
$ws.Range("I102").Value = 842.9
$ws.Range("K102").Value = 842.9
$ws.Range("M102").Value = 779.1

$ws.Range("H122").Value = 2844.303
$ws.Range("I122").Value = 2004.96
$ws.Range("J122").Value = 5467.25
$ws.Range("K122").Value = 6014.88
$ws.Range("L122").Value = 16401.75
$ws.Range("M122").Value = -3564.88
$ws.Range("N122").Value = -21301.75

$ws.Range("H132").Value = 3461072
$ws.Range("I132").Value = 7469.636
$ws.Range("K132").Value = 22408.908
$ws.Range("M132").Value = -19878.908

$ws.Range("H135").Value = 151735.67
$ws.Range("J135").Value = 207606
$ws.Range("L135").Value = 207606
$ws.Range("N135").Value = -217746

$ws.Range("H136").Value = 26769
$ws.Range("I136").Value = 23521.8
$ws.Range("K136").Value = 70565.39999999999
$ws.Range("M136").Value = -68015.39999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 43043.707
$ws.Range("I20").Value = 34779.938
$ws.Range("K20").Value = 34779.938
$ws.Range("M20").Value = -34532.938

$ws.Range("H94").Value = 2314.68
$ws.Range("I94").Value = 1139
$ws.Range("J94").Value = 4813
$ws.Range("K94").Value = 1139
$ws.Range("L94").Value = 4813
$ws.Range("M94").Value = -688
$ws.Range("N94").Value = -5715

$ws.Range("H99").Value = 1301.8572
$ws.Range("I99").Value = 1128.3684
$ws.Range("K99").Value = 1128.3684
$ws.Range("M99").Value = 369.6315999999999

$ws.Range("H105").Value = 2124.0715
$ws.Range("I105").Value = 1979.7693
$ws.Range("J105").Value = 4000
$ws.Range("K105").Value = 1979.7693
$ws.Range("L105").Value = 4000
$ws.Range("M105").Value = -232.7692999999999
$ws.Range("N105").Value = -7494

$ws.Range("H134").Value = 27131.309
$ws.Range("I134").Value = 23670
$ws.Range("K134").Value = 71010
$ws.Range("M134").Value = -68475

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 4253
$ws.Range("I62").Value = 3905
$ws.Range("J62").Value = 4833
$ws.Range("K62").Value = 3905
$ws.Range("L62").Value = 4833
$ws.Range("M62").Value = -3281
$ws.Range("N62").Value = -6081

$ws.Range("H65").Value = 4253
$ws.Range("I65").Value = 3905
$ws.Range("J65").Value = 4833
$ws.Range("K65").Value = 19525
$ws.Range("L65").Value = 24165
$ws.Range("M65").Value = -16405
$ws.Range("N65").Value = -30405

$ws.Range("H107").Value = 393705.1
$ws.Range("I107").Value = 524427.5
$ws.Range("K107").Value = 524427.5
$ws.Range("M107").Value = -522507.5

$ws.Range("H134").Value = 33340536
$ws.Range("I134").Value = 2195.25
$ws.Range("K134").Value = 6585.75
$ws.Range("M134").Value = -4050.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 128.7
$ws.Range("I23").Value = 36
$ws.Range("J23").Value = 151.875
$ws.Range("K23").Value = 108
$ws.Range("L23").Value = 455.625
$ws.Range("M23").Value = 127
$ws.Range("N23").Value = -925.625

$ws.Range("H103").Value = 1767.7222
$ws.Range("I103").Value = 774.2857
$ws.Range("J103").Value = 2399.9092
$ws.Range("K103").Value = 2322.8571
$ws.Range("L103").Value = 7199.7276
$ws.Range("M103").Value = -1443.8571
$ws.Range("N103").Value = -8957.7276

$ws.Range("H117").Value = 1656.9412
$ws.Range("J117").Value = 2339
$ws.Range("L117").Value = 7017
$ws.Range("N117").Value = -13901

$ws.Range("H121").Value = 1444.579
$ws.Range("I121").Value = 466.5
$ws.Range("K121").Value = 1399.5
$ws.Range("M121").Value = -89.5

$ws.Range("H132").Value = 1263.625
$ws.Range("I132").Value = 1307.1818
$ws.Range("J132").Value = 1167.8
$ws.Range("K132").Value = 11764.6362
$ws.Range("L132").Value = 10510.2
$ws.Range("M132").Value = -9234.636200000001
$ws.Range("N132").Value = -15570.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 120.21739
$ws.Range("I2").Value = 134.46666
$ws.Range("K2").Value = 134.46666
$ws.Range("M2").Value = -21.46665999999999

$ws.Range("H69").Value = 47387.145
$ws.Range("J69").Value = 47387.145
$ws.Range("L69").Value = 47387.145
$ws.Range("N69").Value = -48885.145

$ws.Range("H72").Value = 47387.145
$ws.Range("J72").Value = 47387.145
$ws.Range("L72").Value = 142161.435
$ws.Range("N72").Value = -149649.435

$ws.Range("H80").Value = 17393.111
$ws.Range("I80").Value = 11358.4
$ws.Range("J80").Value = 24936.5
$ws.Range("K80").Value = 11358.4
$ws.Range("L80").Value = 24936.5
$ws.Range("M80").Value = -10360.4
$ws.Range("N80").Value = -26932.5

$ws.Range("H83").Value = 17393.111
$ws.Range("I83").Value = 11358.4
$ws.Range("J83").Value = 24936.5
$ws.Range("K83").Value = 56792
$ws.Range("L83").Value = 124682.5
$ws.Range("M83").Value = -51800
$ws.Range("N83").Value = -134666.5

$ws.Range("H97").Value = 957.53125
$ws.Range("I97").Value = 796.5217
$ws.Range("J97").Value = 1369
$ws.Range("K97").Value = 796.5217
$ws.Range("L97").Value = 1369
$ws.Range("M97").Value = -300.5217
$ws.Range("N97").Value = -2361

$ws.Range("H107").Value = 511.90475
$ws.Range("I107").Value = 267.91666
$ws.Range("J107").Value = 837.2222
$ws.Range("K107").Value = 267.91666
$ws.Range("L107").Value = 837.2222
$ws.Range("M107").Value = 1652.08334
$ws.Range("N107").Value = -4677.2222

$ws.Range("H126").Value = 5644.375
$ws.Range("I126").Value = 4582.625
$ws.Range("J126").Value = 6706.125
$ws.Range("K126").Value = 13747.875
$ws.Range("L126").Value = 20118.375
$ws.Range("M126").Value = -11277.875
$ws.Range("N126").Value = -25058.375

$ws.Range("H132").Value = 11720.866
$ws.Range("I132").Value = 4815
$ws.Range("K132").Value = 14445
$ws.Range("M132").Value = -11915

$ws.Range("H134").Value = 17860.6
$ws.Range("J134").Value = 17860.6
$ws.Range("L134").Value = 53581.8
$ws.Range("N134").Value = -58651.8

$ws.Range("H135").Value = 160938.47
$ws.Range("J135").Value = 160938.47
$ws.Range("L135").Value = 160938.47
$ws.Range("N135").Value = -171078.47

$ws.Range("H136").Value = 18677.928
$ws.Range("J136").Value = 18677.928
$ws.Range("L136").Value = 56033.784
$ws.Range("N136").Value = -61133.784

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 324606.03
$ws.Range("I46").Value = 771107.1
$ws.Range("J46").Value = 2133
$ws.Range("K46").Value = 771107.1
$ws.Range("L46").Value = 2133
$ws.Range("M46").Value = -770919.1
$ws.Range("N46").Value = -2509

$ws.Range("H82").Value = 6514.579
$ws.Range("I82").Value = 5398.3335
$ws.Range("J82").Value = 7519.2
$ws.Range("K82").Value = 5398.3335
$ws.Range("L82").Value = 7519.2
$ws.Range("M82").Value = -5037.3335
$ws.Range("N82").Value = -8241.200000000001

$ws.Range("H85").Value = 6514.579
$ws.Range("I85").Value = 5398.3335
$ws.Range("J85").Value = 7519.2
$ws.Range("K85").Value = 5398.3335
$ws.Range("L85").Value = 7519.2
$ws.Range("M85").Value = -4150.3335
$ws.Range("N85").Value = -10015.2

$ws.Range("H93").Value = 9440
$ws.Range("I93").Value = 9788.223
$ws.Range("K93").Value = 9788.223
$ws.Range("M93").Value = -8540.223

$ws.Range("H132").Value = 1150405.9
$ws.Range("I132").Value = 3641.1177
$ws.Range("J132").Value = 2233461.5
$ws.Range("K132").Value = 10923.3531
$ws.Range("L132").Value = 6700384.5
$ws.Range("M132").Value = -8393.3531
$ws.Range("N132").Value = -6705444.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("M25").ClearContents()
$ws.Range("N25").ClearContents()

$ws.Range("H81").Value = 1634.3
$ws.Range("I81").Value = 1260.3334
$ws.Range("J81").Value = 5000
$ws.Range("K81").Value = 2520.6668
$ws.Range("L81").Value = 10000
$ws.Range("M81").Value = -1459.6668
$ws.Range("N81").Value = -12122

$ws.Range("H84").Value = 1634.3
$ws.Range("I84").Value = 1260.3334
$ws.Range("J84").Value = 5000
$ws.Range("K84").Value = 12603.334
$ws.Range("L84").Value = 50000
$ws.Range("M84").Value = -7299.333999999999
$ws.Range("N84").Value = -60608

$ws.Range("H113").Value = 3733.35
$ws.Range("I113").Value = 4485.769
$ws.Range("K113").Value = 13457.307
$ws.Range("M113").Value = -11287.307

$ws.Range("H136").Value = 7744.227
$ws.Range("I136").Value = 2120.4
$ws.Range("J136").Value = 15144
$ws.Range("K136").Value = 6361.200000000001
$ws.Range("L136").Value = 45432
$ws.Range("M136").Value = -3811.200000000001
$ws.Range("N136").Value = -50532
